$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / string updates (safe, non-numeric-looking values)
$ws.Range("D2").Value = "26.919.62"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.550.32"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E8").Value = "  +2.82%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "1.550.61"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "26.912.39"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0₃0706"
$ws.Range("E18").Value = "  +2.95%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("E20").Value = "  +0.48%  "
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("E33").Value = "  +5.38%  "
$ws.Range("D34").Value = "1.411.70"
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("E35").Value = "  +3.03%  "
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  +1.19%  "
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("E42").Value = "  +3.40%  "
$ws.Range("E43").Value = "  +2.86%  "
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("E45").Value = "  +1.56%  "
$ws.Range("D47").Value = "1.685.23"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").Value = "0.0₇0999"
$ws.Range("E50").Value = "  +3.09%  "
$ws.Range("E51").Value = "  +0.83%  "

# Values that look numeric must be forced to text so Excel does not
# reinterpret them (e.g. "61.67" -> 61.67 as a float). We temporarily
# apply a text number format, assign the value, then clear the format
# override so the cell keeps its original (default) style.
$numericLookingCells = @{
    "D5" = "206.25"
    "D8" = "22.10"
    "D10" = "0.0589"
    "D11" = "0.0856"
    "D17" = "61.67"
    "D19" = "217.39"
    "D20" = "7.27"
    "D23" = "9.23"
    "D25" = "153.91"
    "D27" = "14.94"
    "D33" = "3.12"
    "D38" = "0.0165"
    "D39" = "0.528"
    "D44" = "0.998"
    "D45" = "64.64"
    "D46" = "1.74"
    "D48" = "87.51"
    "D49" = "0.0516"
    "D51" = "0.0961"
}
foreach ($ref in $numericLookingCells.Keys) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $numericLookingCells[$ref]
    $r.ClearFormats()
}
